$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "BoxOffice" column (column AB / 28),
# shifting BoxOffice and everything after it one column to the right.
$ws.Columns.Item(28).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 28).Value = "imdbRating"

# IMDB rating values for rows 2..134 (one per movie row), in order.
$ratings = @(7.2,7.1,6.3,5.3,6.1,5.8,6.2,6.3,6.6,6.6,5.9,7.2,6.6,7.6,7.5,4.3,6.3,6.4,6.4,8.1,6.2,6.8,6.9,7,3.9,8,6.3,8.1,6.7,6,6,4.9,5.6,5.2,6,5.5,6.6,6.3,7.5,6.8,5.7,4.6,6.8,6.4,7.2,5.1,6.5,6.6,6.5,6.4,6.3,6,6,7,6.1,7.5,6.7,6.2,6.2,6.9,6.8,5.9,6.7,7.8,8.5,7.1,6.1,6.6,5.8,3.9,6.8,8,6.9,6.4,6.1,5.3,6,6.2,6.8,6.6,6.2,6.2,5.8,7.4,5.4,7,6.9,5.9,6.7,3,7.8,8,7.2,6,7.7,8,6.5,7.8,6.7,5.5,5.4,6.9,6.7,5.3,4.4,6.4,6.2,6.3,7.9,6.4,7.4,5,6,7.1,7.3,7.8,6.4,6.1,8.300000000000001,6.4,5.8,7.4,7.8,8.300000000000001,7.1,6.1,6.1,6.9,5.1,7.4,6.5,6.5,7.2)

for ($i = 0; $i -lt $ratings.Length; $i++) {
    $ws.Cells.Item($i + 2, 28).Value = $ratings[$i]
}
